$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain plain text so values like "313.91" or "0.4860" are preserved exactly
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.588.08"
$ws.Range("E2").Value = "  +2.47%  "
$ws.Range("D3").Value = "1.878.86"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("D4").Value = "1.016"
$ws.Range("E4").Value = "  +0.88%  "
$ws.Range("D5").Value = "313.91"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("D7").Value = "0.4804"
$ws.Range("E7").Value = "  +1.49%  "
$ws.Range("D8").Value = "0.3790"
$ws.Range("E8").Value = "  +2.96%  "
$ws.Range("D9").Value = "0.07399"
$ws.Range("E9").Value = "  +2.42%  "
$ws.Range("D10").Value = "0.9432"
$ws.Range("D11").Value = "20.74"
$ws.Range("E11").Value = "  +5.66%  "
$ws.Range("D12").Value = "0.07877"
$ws.Range("E12").Value = "  +3.32%  "
$ws.Range("D13").Value = "1.886.86"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").Value = "5.459"
$ws.Range("E14").Value = "  +2.75%  "
$ws.Range("D15").Value = "6.618"
$ws.Range("E15").Value = "  +3.46%  "
$ws.Range("D16").Value = "91.41"
$ws.Range("E16").Value = "  +3.47%  "
$ws.Range("D17").Value = "1.016"
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("D18").Value = "0.000009016"
$ws.Range("E18").Value = "  +3.97%  "
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").Value = "14.99"
$ws.Range("E20").Value = "  +2.88%  "
$ws.Range("D21").Value = "27.594.09"
$ws.Range("E21").Value = "  +2.39%  "
$ws.Range("D22").Value = "5.149"
$ws.Range("E22").Value = "  +2.20%  "
$ws.Range("D23").Value = "10.79"
$ws.Range("E23").Value = "  +1.10%  "
$ws.Range("E24").Value = "  +2.69%  "
$ws.Range("D25").Value = "153.72"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("D26").Value = "18.61"
$ws.Range("E26").Value = "  +2.44%  "
$ws.Range("D27").Value = "2.032"
$ws.Range("E27").Value = "  +1.49%  "
$ws.Range("D28").Value = "116.31"
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("D29").Value = "5.015"
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("D30").Value = "0.08937"
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("D31").Value = "3.327"
$ws.Range("E31").Value = "  +1.34%  "
$ws.Range("D32").Value = "1.216"
$ws.Range("E32").Value = "  +4.22%  "
$ws.Range("E33").Value = "  +2.80%  "
$ws.Range("D34").Value = "0.7532"
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("D35").Value = "2.700"
$ws.Range("E35").Value = "  -2.41%  "
$ws.Range("D36").Value = "0.02081"
$ws.Range("E36").Value = "  +6.67%  "
$ws.Range("D37").Value = "1.122"
$ws.Range("E37").Value = "  +2.94%  "
$ws.Range("D38").Value = "0.05315"
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("D39").Value = "3.014"
$ws.Range("E39").Value = "  +1.64%  "
$ws.Range("D40").Value = "0.5384"
$ws.Range("E40").Value = "  +3.43%  "
$ws.Range("D41").Value = "7.124"
$ws.Range("E41").Value = "  +2.99%  "
$ws.Range("D42").Value = "0.1526"
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("D43").Value = "8.469"
$ws.Range("E43").Value = "  +3.21%  "

# Rows 44 and 45 swap Decentraland/EnergySwap entries (with updated price/volume)
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "0.4860"
$ws.Range("E44").Value = "  +3.44%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "10.64"
$ws.Range("E45").Value = "  +0.70%  "

$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("D47").Value = "1.668"
$ws.Range("E47").Value = "  +4.18%  "
$ws.Range("D48").Value = "103.23"
$ws.Range("E48").Value = "  +1.18%  "
$ws.Range("D49").Value = "67.42"
$ws.Range("E49").Value = "  +2.95%  "
$ws.Range("D50").Value = "0.06108"
$ws.Range("E50").Value = "  +1.22%  "
$ws.Range("D51").Value = "0.9029"
$ws.Range("E51").Value = "  +1.99%  "
